# Update Work Week and Social Spending
# (Commit message retained verbatim from source repo; the actual edit
# refreshes the Nicaragua GDP-per-Capita "Data" series for years 1920-2008
# and appends new rows for years 2009-2016.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# row, year, value (value kept as text to match the source workbook's
# shared-string-typed "Data" column)
$rowsData = @(
    @(2, 1920, "2007"),
    @(3, 1921, "2083"),
    @(4, 1922, "1878"),
    @(5, 1923, "2010"),
    @(6, 1924, "2098"),
    @(7, 1925, "2316"),
    @(8, 1926, "1983"),
    @(9, 1927, "1992"),
    @(10, 1928, "2523"),
    @(11, 1929, "2778"),
    @(12, 1930, "2246"),
    @(13, 1931, "2071"),
    @(14, 1932, "1863"),
    @(15, 1933, "2311"),
    @(16, 1934, "2067"),
    @(17, 1935, "2043"),
    @(18, 1936, "1583"),
    @(19, 1937, "1672"),
    @(20, 1938, "1707"),
    @(21, 1939, "2042"),
    @(22, 1940, "2179"),
    @(23, 1941, "2354"),
    @(24, 1942, "2214"),
    @(25, 1943, "2373"),
    @(26, 1944, "2299"),
    @(27, 1945, "2259"),
    @(28, 1946, "2377"),
    @(29, 1947, "2310"),
    @(30, 1948, "2460"),
    @(31, 1949, "2345"),
    @(32, 1950, "2565"),
    @(33, 1951, "2657"),
    @(34, 1952, "3016"),
    @(35, 1953, "2997"),
    @(36, 1954, "3178"),
    @(37, 1955, "3288"),
    @(38, 1956, "3186"),
    @(39, 1957, "3351"),
    @(40, 1958, "3258"),
    @(41, 1959, "3204"),
    @(42, 1960, "3148"),
    @(43, 1961, "3277"),
    @(44, 1962, "3521"),
    @(45, 1963, "3781"),
    @(46, 1964, "4092"),
    @(47, 1965, "4340"),
    @(48, 1966, "4342"),
    @(49, 1967, "4500"),
    @(50, 1968, "4419"),
    @(51, 1969, "4564"),
    @(52, 1970, "4463"),
    @(53, 1971, "4533"),
    @(54, 1972, "4543"),
    @(55, 1973, "4637"),
    @(56, 1974, "5064"),
    @(57, 1975, "5013"),
    @(58, 1976, "5099"),
    @(59, 1977, "5247"),
    @(60, 1978, "4771"),
    @(61, 1979, "3405"),
    @(62, 1980, "3414"),
    @(63, 1981, "3480"),
    @(64, 1982, "3360"),
    @(65, 1983, "3437"),
    @(66, 1984, "3308"),
    @(67, 1985, "3107"),
    @(68, 1986, "3011"),
    @(69, 1987, "2922"),
    @(70, 1988, "2498"),
    @(71, 1989, "2393"),
    @(72, 1990, "2291"),
    @(73, 1991, "2209.79424830776"),
    @(74, 1992, "2163.41659168443"),
    @(75, 1993, "2114.69189753674"),
    @(76, 1994, "2154.29901195949"),
    @(77, 1995, "2264.09251825208"),
    @(78, 1996, "2395.20217328583"),
    @(79, 1997, "2482.69082400451"),
    @(80, 1998, "2569.82833048217"),
    @(81, 1999, "2748.12700024576"),
    @(82, 2000, "2858.5269124027"),
    @(83, 2001, "2944.24452676668"),
    @(84, 2002, "2969.4575848692"),
    @(85, 2003, "3049.14999132693"),
    @(86, 2004, "3220.26103652372"),
    @(87, 2005, "3370.77750867683"),
    @(88, 2006, "3525.553651962"),
    @(89, 2007, "3729.85143334289"),
    @(90, 2008, "3856.55228301069"),
    @(91, 2009, "3771.94704872872"),
    @(92, 2010, "3916.96592569171"),
    @(93, 2011, "4189"),
    @(94, 2012, "4376"),
    @(95, 2013, "4527"),
    @(96, 2014, "4685"),
    @(97, 2015, "4866"),
    @(98, 2016, "5045")
)

foreach ($entry in $rowsData) {
    $r = $entry[0]
    $year = $entry[1]
    $val = $entry[2]

    $ws.Cells.Item($r, 1).Value = 558
    $ws.Cells.Item($r, 2).Value = "Nicaragua"
    $ws.Cells.Item($r, 3).Value = "GDP per Capita"
    $ws.Cells.Item($r, 4).Value = $year

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $val
}
